$d = $word.ActiveDocument

function Insert-BodyXml($bodyXml) {
    $last = $d.Paragraphs.Last
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $last.Range.InsertXML($pkg)
}

$newParas = '<w:p><w:r><w:t xml:space="preserve">Contains the genetic information of the prokaryotic cell. The nucleoid determines how the cell will develop and grow. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Responsible for </w:t></w:r><w:r><w:t xml:space="preserve">making proteins. These proteins </w:t></w:r><w:r><w:t>are what keeps the cell alive by providing energy for movement, repair, and growth.</w:t></w:r></w:p><w:p><w:r><w:t>These are extra DNA molecules that are commonly found in bacteria, and sometimes in other cells. The</w:t></w:r><w:r><w:t>y allow cells to mutate and adapt to an ever-changing environment.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">An affinity for extremely hot temperature. These </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>archaeans</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">can thrive in environments near volcanos, hot springs, and acidic soils. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">These </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>archaeans</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> release methane as a result of digesting materials such as carbon and </w:t></w:r><w:r><w:t>hydrogen</w:t></w:r><w:r><w:t xml:space="preserve">. They are known to play a role in breaking up </w:t></w:r><w:r><w:t>materials for other cells to consume.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">An affinity for salt. These </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>archaeans</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>can withstand the effects of salt which causes dehydration. They are also known to be resistant to UV radiation, giving them a reddish look.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:t>single flagellum</w:t></w:r><w:r><w:t xml:space="preserve"> that favors long distance travel.</w:t></w:r></w:p><w:p><w:r><w:t>Multiple flagella in one polar end of the cell that favors frequent twists and turns.</w:t></w:r></w:p><w:p><w:r><w:t>Multiple flagella around the cell that helps with even more frequent twists and turns.</w:t></w:r></w:p><w:p><w:r><w:t>A prokaryote that metabolizes methane. This particular type is anaerobic, which means it can live with little to no oxygen.</w:t></w:r></w:p><w:p><w:r><w:t>Capable of harnessing energy from sunlight to produce food from inorganic matter. Its chlorophyl component gives it a green color.</w:t></w:r></w:p><w:p><w:r><w:t>Eats anything organic that contains carbon and hydrogen, such as glucose.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">This particular bacteria feeds from </w:t></w:r><w:r><w:t>within</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>the</w:t></w:r><w:r><w:t xml:space="preserve"> host, draining it of all its energy. Once engulfed</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>it secretes special proteins that will prevent itself from being consumed.</w:t></w:r></w:p>'

Insert-BodyXml($newParas)

Write-Output $d.Paragraphs.Count
